# Update the two-digit/one-digit division worksheet numbers.
# Each "old" cell text is unique in the document, so a targeted
# Find/Replace (whole document, match case, no wildcards) for each
# pair is safe. The pairs are applied in the same order they appear
# in the document so that a later replacement's "new" text never
# collides with an earlier replacement's "old" text search.

$d = $word.ActiveDocument

$pairs = @(
    @("50÷5=", "68÷6="),
    @("43÷8=", "95÷5="),
    @("36÷2=", "42÷7="),
    @("51÷6=", "75÷2="),
    @("78÷3=", "83÷8="),
    @("96÷9=", "35÷4="),
    @("90÷6=", "41÷7="),
    @("11÷7=", "60÷5="),
    @("25÷6=", "31÷8="),
    @("52÷6=", "21÷4="),
    @("13÷7=", "16÷3="),
    @("93÷7=", "29÷5="),
    @("55÷7=", "98÷7="),
    @("15÷4=", "72÷9="),
    @("84÷5=", "84÷7="),
    @("43÷9=", "62÷8="),
    @("29÷2=", "74÷5="),
    @("59÷6=", "60÷7="),
    @("48÷3=", "26÷4="),
    @("17÷2=", "93÷7="),
    @("89÷4=", "67÷2="),
    @("85÷8=", "90÷6="),
    @("12÷8=", "51÷9="),
    @("84÷8=", "78÷9="),
    @("36÷5=", "15÷5=")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
}
